$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E, shifting old E (aiAnalysis) to F and old F (ResultLink) to G
$ws.Columns("E:E").Insert()

$values = @{}
$values['E1'] = 'Test_Submitted_Time'
$values['F1'] = 'aiAnalysis'
$values['G1'] = 'ResultLink'
$values['E2'] = '2025-01-09 | 05:10:44 PM'
$values['F2'] = 'No solution is fetched'
$values['G2'] = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2FPz1cTonXEa5S0eGjjClwky18L7ZLivoOknU9lnguXeJnyFo9x8IKGQYVqJ%2BeXXRYNM29wXlqjOFR6o1N3rli9xDNMkL5K3aS4kQ0TGWrQxdZDKotAc25sNzu4mYKKKe073uRDcl84QA%3D%3D'
$values['E3'] = '2025-01-09 | 05:56:18 PM'
$values['F3'] = 'No solution is fetched'
$values['G3'] = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX19Kq5F7MGU0pdxj8vGZecI%2Bf6A25%2BE2a2a3cPUPOX2ucPTgSuvn4LaymTQbC%2FGHGeMYCqLp7FNC9DA%2FqK2MhC83XLqizF5blm8W8Xz1JDu0%2FeOOj4HucLvQj7Ykt7%2FUz40sqgHp5cxvgA%3D%3D'
$values['E4'] = '2025-01-09 | 06:00:36 PM'
$values['F4'] = 'Here''s a line-by-line analysis of the provided solution:
In the `Vehicle` class, the properties are declared but not set. This might lead to incorrect data handling.
In the `Program` class, the `AddVehicle` method is not actually adding the vehicle to the database. It only creates a `DataTable` and sets it to `ds`. The code should use an `Insert` statement to add the vehicle to the database.
The `ListVehicles` method in the `Program` class is also not actually listing the vehicles. It creates a `DataSet` and a `DataAdapter`, but does not traverse the data tables to display the vehicle information.
The `DeleteVehicle` method is declared but not implemented. It does not have a logic to delete a vehicle from the database.
The connection string is hard-coded, which might not be suitable in a real-world application.
The `success` class is an empty class, and it seems like a leftover from a previous version of the code.
Final Analysis: The provided solution has several logical and runtime errors. The `AddVehicle`, `ListVehicles`, and `DeleteVehicle` methods are not implemented correctly, which will lead to incorrect data handling and potential database errors. The code also has issues with data security, as the connection string is hard-coded.'
$values['G4'] = ''
$values['E5'] = '2025-01-09 | 05:38:21 PM'
$values['F5'] = 'Based on the provided solution, the code has several errors and warnings. 
1. The `Vehicle` class is not properly designed. The properties `VehicleID`, `Make`, `Model`, `Year`, and `LastMaintenanceDate` are not private and are publicly accessible, which is a potential security risk.
2. The `Program` class is trying to create a new instance of `Vehicle` in the `Main` method, which is not necessary. The `Main` method can simply get user input and call the corresponding methods based on the user''s choice.
3. The `AddVehicle` method is inserting values into the `Vehicles` table, but the `VehicleID` property of the `Vehicle` class is not set. This will result in a primary key constraint error when trying to insert the vehicle.
4. The `ListVehicles` method is trying to execute the SQL query, but it''s not populating any variables to print the results.
5. The `DeleteVehicle` method is trying to delete a vehicle from the `Vehicles` table, but it''s not handling the case where the vehicle does not exist.
6. There is no effort to handle exceptions that may occur during database operations.
7. The code does not have any user-friendly error handling. For example, if the user enters an invalid choice, the program simply prints "Invalid choice" and terminates.
8. The methods `AddVehicle`, `ListVehicles`, and `DeleteVehicle` are all static methods, but they are not declared as such. 
In conclusion, the code is not properly designed and has several errors and warnings. It needs to be rewritten to follow best practices and handle exceptions properly.'
$values['G5'] = ''
$values['E6'] = '2025-01-09 | 05:44:23 PM'
$values['F6'] = 'Here''s the analysis of the provided solution:
* The `Vehicle` class is missing in the solution, which is required according to the problem statement.
* The `Program` class contains several methods, including `AddVehicle`, `ListVehicles`, and `DeleteVehicle`, but they are not properly implemented. For example, in `AddVehicle`, the `DataColumn` names are not specified, and in `ListVehicles` and `DeleteVehicle`, the `SqlDataAdapter` and `SqlCommandBuilder` are not used correctly.
* The code also contains several issues with data typing, such as using `DataRow` without specifying the data type, and using `string` instead of `int` for vehicle ID.
Overall, the solution lacks proper implementation of classes, methods, and data types, resulting in compilation errors and runtime errors. The code also fails to follow the problem statement requirements, such as creating a Vehicle class and implementing CRUD operations correctly.
Final Analysis: The provided solution fails to meet the requirements of the problem statement due to inadequate implementation of classes, methods, and data types, resulting in compilation errors and runtime errors. The code lacks proper use of `SqlDataAdapter` and `SqlCommandBuilder` and fails to follow the required data typing conventions.'
$values['G6'] = ''
$values['E7'] = '2025-01-09 | 05:46:52 PM'
$values['F7'] = ' Failure Analysis:
1. In the "Vehicle.cs" file, the Vehicle class is not correctly defined. The class name is misspelled as "vehicle" instead of "Vehicle".
Failure Code: `public class vehicle`
Correct Code: `public class Vehicle`
2. In the "Program.cs" file, the "AddVehicle" method has no implementation. It is empty and does not perform any action.
Failure Code: `public static void AddVehicle(vehicle vehicle) { }`
Correct Code: `public static void AddVehicle(Vehicle vehicle) { // Insert vehicle into database }`
3. In the "Program.cs" file, the "ListVehicles" method has no implementation. It is empty and does not perform any action.
Failure Code: `public static void ListVehicles() { }`
Correct Code: `public static void ListVehicles() { // Retrieve and display all vehicles }`
4. In the "Program.cs" file, the "DeleteVehicle" method has no implementation. It is empty and does not perform any action.
Failure Code: `public static void DeleteVehicle(int VehicleId) { }`
Correct Code: `public static void DeleteVehicle(int VehicleId) { // Delete the vehicle from the database }`
Final Analysis:
The provided solution has several syntax and logical errors. The Vehicle class is not correctly defined, and the methods in the Program class lack implementation. The error analysis highlights the specific lines of code that require correction.'
$values['G7'] = ''
$values['E8'] = '2025-01-09 | 05:58:40 PM'
$values['F8'] = 'The provided solution does not contain a fully functional implementation of the Vehicle Management System as described. The analyzed files are as follows:
1. Vehicle.cs: This file does not contain the Vehicle class with properties as described. It is missing the crucial implementation of the class with properties like VehicleID, Make, Model, Year, and LastMaintenanceDate.
Failure analysis: The Vehicle class is missing its implementation, and its properties are not defined, which is a critical component of the system.
2. Program.cs: This file has some basic setup code for the program, but it does not contain the implementation of the AddVehicle, ListVehicles, DeleteVehicle, or the Main method as described.
Failure analysis: The code is incomplete, and the methods that should be implemented are missing, which means the program cannot be run as described.
Final analysis: The provided solution is incomplete and does not implement the Vehicle class or the required methods, making it impossible to run or test the system.'
$values['G8'] = ''
$values['E9'] = '2025-01-09 | 05:16:19 PM'
$values['F9'] = 'The provided solution has several issues. In the `Vehicle` class, the property names are not camel-cased, which is against the conventional naming standards in C#. Additionally, the `ConnectionString` property in the `ConnectionStringProvider` class is not static, which may cause issues when trying to access the connection string.
In the `AddVehicle` method, the code is attempting to insert a new vehicle into the database using `SqlDataAdapter`. However, this method is not designed for inserting data and may cause runtime errors. The correct approach is to use `SqlCommand` to execute an INSERT statement.
The final analysis is that the provided solution has several syntax and logical errors that would prevent it from compiling or running correctly. It also lacks proper error handling and connection string management, which could lead to runtime errors and security issues.'
$values['G9'] = ''
$values['E10'] = '2025-01-09 | 05:59:33 PM'
$values['F10'] = 'No solution is fetched'
$values['G10'] = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2FBHtt1n3Z0Kxww9xueXMAIQsoaslwp7ocuGeh6LxE3vutsv3JbfHNCyu9tzcOyIuQUqZQKvQxYvJKk05uKwByHYCPXrH5FzyRDfZPdc4R0wyAe49O6QROrH2ERX8fYwayubOBz%2BGtvkg%3D%3D'
$values['E11'] = '2025-01-09 | 05:44:45 PM'
$values['F11'] = 'The provided solution file, Program.cs, has several issues that need to be addressed.
Firstly, the Vehicle class is defined with properties that do not match the requirements specified in the problem statement. The properties should be public, but they are declared as auto-implemented properties, which is not correct.
Secondly, the AddVehicle method is not implementing the required functionality. It creates a SqlDataAdapter but does not use it to insert the new vehicle into the database.
Thirdly, the ListVehicles method is not implemented correctly. It creates a SqlDataAdapter but does not use it to retrieve the list of vehicles from the database.
Fourthly, the DeleteVehicle method is not implemented correctly. It creates a SqlDataAdapter but does not use it to delete the vehicle with the specified ID from the database.
Final Analysis: The provided solution file has several logical and syntax errors that need to be addressed. The Vehicle class properties need to be declared as public, the AddVehicle method needs to use the SqlDataAdapter to insert the new vehicle into the database, the ListVehicles method needs to use the SqlDataAdapter to retrieve the list of vehicles from the database, and the DeleteVehicle method needs to use the SqlDataAdapter to delete the vehicle with the specified ID from the database.'
$values['G11'] = ''

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
